# Revert the earlier "edited excel files" commit: strip the "/24" CIDR
# suffix back off the two IP-address cells on Sheet1, and restore the
# last-saved selection to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# F2 (SW1 IP address) and F3 (SW2 IP address) were "192.168.10.5/24" and
# "10.1.1.5/24" - put them back to plain IPs without the prefix length.
$ws.Range("F2").Value = "192.168.10.5"
$ws.Range("F3").Value = "10.1.1.5"

# Restore the saved cursor/selection position to D5.
$ws.Range("D5").Select()
